# Survey workbook template update:
# - Replace the "A3/A4/A5" free-text answer columns with structured
#   Solution / StartSet / Difficulty / Slope columns used by the new
#   backend file scheme.
# - Re-point the active selection and tidy up the "StartSet" column width.
# - Apply a basic page setup (paper size / orientation) for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey")

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("E1").Value = "Solution"
$ws.Range("F1").Value = "StartSet"
$ws.Range("G1").Value = "Difficulty"
$ws.Range("H1").Value = "Slope"

# --- Row 2 (Pizza?) --------------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "X"
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()

# --- Row 3 (Döner?) --------------------------------------------------------
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "X"
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()

# --- Row 4 (Nudeln?) -------------------------------------------------------
$ws.Range("E4").Value = 1
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# --- Column width for the now wider "StartSet" header ----------------------
$ws.Range("F1").ColumnWidth = 12.14

# --- Selection / active cell moves from H4 to F2 ---------------------------
[void]$ws.Range("F2").Select()

# --- Basic print page setup -------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
